$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (d=9 row, A7 label)
$ws.Range("B7").Value = 97.73187545085479
$ws.Range("C7").Value = 97.74825703922609
$ws.Range("D7").Value = 97.79188635546379
$ws.Range("E7").Value = 97.72524485738487

# Row 8 (d=10 row, A8 label)
$ws.Range("B8").Value = 97.27415082650515
$ws.Range("C8").Value = 97.22816963485333
$ws.Range("D8").Value = 97.28024873775524
$ws.Range("E8").Value = 97.2303666519815

# Row 9 (d=11 row, A9 label)
$ws.Range("B9").Value = 95.90003212440614
$ws.Range("C9").Value = 95.89506949893246
$ws.Range("D9").Value = 95.88467961110342
$ws.Range("E9").Value = 95.92727481418855
